$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "data"

# Move the "omzet" column (column B) so it becomes the last column,
# shifting "prijs" and "promotiekosten" one position to the left.
$ws.Columns.Item(2).Cut() | Out-Null
$ws.Columns.Item(5).Insert() | Out-Null

# The cut/insert reorders the underlying cell data (and therefore the
# table header text) correctly, but the ListObject's cached column
# metadata needs to be nudged so it re-derives column names from the
# (now reordered) header cells.
$ws.Range("B1").Value = $ws.Range("B1").Value()
$ws.Range("C1").Value = $ws.Range("C1").Value()
$ws.Range("D1").Value = $ws.Range("D1").Value()

# Rename the table itself
$lo = $ws.ListObjects.Item(1)
$lo.Name = "TabelEnergierepen"

# Update page setup (paper size / orientation)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Update the active cell selection
$ws.Range("A5").Select() | Out-Null
